$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format so numeric-looking strings
# (e.g. "247.37", "0.9996") are preserved as text instead of being
# coerced to numbers, matching the original inline-string cell type.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '30.517.74'
$ws.Range("D3").Value = '1.880.95'
$ws.Range("D4").Value = '0.9996'
$ws.Range("D5").Value = '247.37'
$ws.Range("D6").Value = '0.9994'
$ws.Range("D7").Value = '0.4754'
$ws.Range("D8").Value = '0.2903'
$ws.Range("D9").Value = '0.06527'
$ws.Range("D10").Value = '22.02'
$ws.Range("D11").Value = '0.07734'
$ws.Range("D12").Value = '97.60'
$ws.Range("D13").Value = '0.7445'
$ws.Range("D14").Value = '1.877.11'
$ws.Range("D15").Value = '5.144'
$ws.Range("D16").Value = '275.50'
$ws.Range("D17").Value = '30.493.69'
$ws.Range("D18").Value = '13.62'
$ws.Range("D19").Value = '0.000007573'
$ws.Range("D20").Value = '0.9998'
$ws.Range("D21").Value = '2.124.92'
$ws.Range("D23").Value = '5.275'
$ws.Range("D24").Value = '6.197'
$ws.Range("D25").Value = '9.339'
$ws.Range("D26").Value = '163.38'
$ws.Range("D27").Value = '18.92'
$ws.Range("D29").Value = '1.371'
$ws.Range("D30").Value = '0.09967'
$ws.Range("D31").Value = '1.522'
$ws.Range("D32").Value = '4.323'
$ws.Range("D33").Value = '4.075'
$ws.Range("D34").Value = '0.04800'
$ws.Range("D35").Value = '1.130'
$ws.Range("D36").Value = '0.7021'
$ws.Range("D37").Value = '2.713'
$ws.Range("D38").Value = '0.01877'
$ws.Range("D39").Value = '2.732'
$ws.Range("D40").Value = '6.338'
$ws.Range("D41").Value = '1.964'
$ws.Range("D42").Value = '71.24'
$ws.Range("D43").Value = '0.4255'
$ws.Range("D44").Value = '0.8406'
$ws.Range("D45").Value = '0.9994'
$ws.Range("D46").Value = '102.88'
$ws.Range("D47").Value = '9.257'
$ws.Range("D48").Value = '7.104'
$ws.Range("D49").Value = '35.61'
$ws.Range("D50").Value = '923.75'
$ws.Range("D51").Value = '0.3901'

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("E5").Value = '  +5.80%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("E7").Value = '  +1.52%  '
$ws.Range("E8").Value = '  +2.67%  '
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("E10").Value = '  +4.97%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("E12").Value = '  +4.44%  '
$ws.Range("E13").Value = '  +9.86%  '
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("E16").Value = '  +3.43%  '
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("E18").Value = '  +2.51%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  +2.95%  '
$ws.Range("E24").Value = '  +1.74%  '
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("E28").Value = '  +3.72%  '
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("E31").Value = '  +5.18%  '
$ws.Range("E32").Value = '  +3.17%  '
$ws.Range("E33").Value = '  +2.63%  '
$ws.Range("E34").Value = '  +3.46%  '
$ws.Range("E35").Value = '  +1.64%  '
$ws.Range("E36").Value = '  +2.26%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +2.78%  '
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("E41").Value = '  +4.99%  '
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("E43").Value = '  +5.40%  '
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("E49").Value = '  +4.95%  '
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("E51").Value = '  +4.59%  '

# --- Row 51: coin renamed from Cronos to Decentraland ---
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

# Clear the temporary text-number-format so the cell style index
# reverts to the default/original (no explicit style), matching
# the source workbook where these cells carry no "s" attribute.
$ws.Range("D2:E51").Style = "Normal"
